$wb = $excel.ActiveWorkbook

# ALC row 5
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 246.4
$ws.Range("I5").Value = 239.4
$ws.Range("K5").Value = 239.4
$ws.Range("M5").Value = -124.4

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 65443.516
$ws.Range("I17").Value = 63
$ws.Range("J17").Value = 66880.45
$ws.Range("K17").Value = 189
$ws.Range("L17").Value = 200641.35
$ws.Range("M17").Value = -21
$ws.Range("N17").Value = -200977.35

# ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 692.1429000000001
$ws.Range("I18").Value = 692.1429000000001
$ws.Range("K18").Value = 692.1429000000001
$ws.Range("M18").Value = -408.1429000000001

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1787.2727
$ws.Range("I28").Value = 1553
$ws.Range("K28").Value = 1553
$ws.Range("M28").Value = -1068

# ALC row 59
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H59").Value = 2049.8333
$ws.Range("J59").Value = 2049.8333
$ws.Range("L59").Value = 6149.499899999999
$ws.Range("N59").Value = -7263.499899999999

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5984.077
$ws.Range("I62").Value = 3757.8
$ws.Range("J62").Value = 7375.5
$ws.Range("K62").Value = 3757.8
$ws.Range("L62").Value = 7375.5
$ws.Range("M62").Value = -3133.8
$ws.Range("N62").Value = -8623.5

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 5984.077
$ws.Range("I65").Value = 3757.8
$ws.Range("J65").Value = 7375.5
$ws.Range("K65").Value = 18789
$ws.Range("L65").Value = 36877.5
$ws.Range("M65").Value = -15669
$ws.Range("N65").Value = -43117.5

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 8574.333000000001
$ws.Range("I86").Value = 7475
$ws.Range("J86").Value = 9124
$ws.Range("K86").Value = 7475
$ws.Range("L86").Value = 9124
$ws.Range("M86").Value = -6352
$ws.Range("N86").Value = -11370

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 8574.333000000001
$ws.Range("I89").Value = 7475
$ws.Range("J89").Value = 9124
$ws.Range("K89").Value = 37375
$ws.Range("L89").Value = 45620
$ws.Range("M89").Value = -31759
$ws.Range("N89").Value = -56852

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1032.1428
$ws.Range("I112").Value = 393.33334
$ws.Range("J112").Value = 1206.3636
$ws.Range("K112").Value = 1180.00002
$ws.Range("L112").Value = 3619.0908
$ws.Range("M112").Value = -72.00001999999995
$ws.Range("N112").Value = -5835.0908

# ALC row 123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 23333.334
$ws.Range("J123").Value = 23333.334
$ws.Range("L123").Value = 23333.334
$ws.Range("N123").Value = -33133.334

# ALC row 126
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H126").Value = 24000
$ws.Range("J126").Value = 24000
$ws.Range("L126").Value = 24000
$ws.Range("N126").Value = -33880

# ALC row 130
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 25555.555
$ws.Range("J130").Value = 25555.555
$ws.Range("L130").Value = 25555.555
$ws.Range("N130").Value = -35595.555

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4241.8965
$ws.Range("I132").Value = 2334.389
$ws.Range("J132").Value = 7363.273
$ws.Range("K132").Value = 7003.167
$ws.Range("L132").Value = 22089.819
$ws.Range("M132").Value = -4473.167
$ws.Range("N132").Value = -27149.819

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4002225.2
$ws.Range("I138").Value = 1242.2703
$ws.Range("J138").Value = 15389638
$ws.Range("K138").Value = 3726.810899999999
$ws.Range("L138").Value = 46168914
$ws.Range("M138").Value = 1413.189100000001
$ws.Range("N138").Value = -46179194

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 109.75
$ws.Range("I5").Value = 109.75
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 109.75
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = 2.25

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3115.8667
$ws.Range("I32").Value = 2473.768
$ws.Range("J32").Value = 10500
$ws.Range("K32").Value = 2473.768
$ws.Range("L32").Value = 10500
$ws.Range("M32").Value = -2186.768
$ws.Range("N32").Value = -11074

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1998.9231
$ws.Range("I45").Value = 1963.5555
$ws.Range("J45").Value = 2078.5
$ws.Range("K45").Value = 1963.5555
$ws.Range("L45").Value = 2078.5
$ws.Range("M45").Value = -1586.5555
$ws.Range("N45").Value = -2832.5

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 862.38464
$ws.Range("I97").Value = 850.55
$ws.Range("J97").Value = 901.8333
$ws.Range("K97").Value = 850.55
$ws.Range("L97").Value = 901.8333
$ws.Range("M97").Value = -354.55
$ws.Range("N97").Value = -1893.8333

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 109.75
$ws.Range("I4").Value = 109.75
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 109.75
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = 5.25

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 546.5454999999999
$ws.Range("I94").Value = 458.375
$ws.Range("J94").Value = 781.6667
$ws.Range("K94").Value = 458.375
$ws.Range("L94").Value = 781.6667
$ws.Range("M94").Value = -7.375
$ws.Range("N94").Value = -1683.6667

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 68183860
$ws.Range("I31").Value = 83334860
$ws.Range("J31").Value = 50002650
$ws.Range("K31").Value = 83334860
$ws.Range("L31").Value = 50002650
$ws.Range("M31").Value = -83334565
$ws.Range("N31").Value = -50003240

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 68183860
$ws.Range("I34").Value = 83334860
$ws.Range("J34").Value = 50002650
$ws.Range("K34").Value = 83334860
$ws.Range("L34").Value = 50002650
$ws.Range("M34").Value = -83334658
$ws.Range("N34").Value = -50003054

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2516.0527
$ws.Range("I62").Value = 2214.6428
$ws.Range("K62").Value = 2214.6428
$ws.Range("M62").Value = -1590.6428

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2516.0527
$ws.Range("I65").Value = 2214.6428
$ws.Range("K65").Value = 11073.214
$ws.Range("M65").Value = -7953.214

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1538.4166
$ws.Range("I107").Value = 503.14285
$ws.Range("K107").Value = 503.14285
$ws.Range("M107").Value = 1416.85715

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 753.64703
$ws.Range("J5").Value = 1400
$ws.Range("L5").Value = 4200
$ws.Range("N5").Value = -4424

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1136.6
$ws.Range("J122").Value = 1740.4615
$ws.Range("L122").Value = 15664.1535
$ws.Range("N122").Value = -20564.1535

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 872.25
$ws.Range("J131").Value = 890.15466
$ws.Range("L131").Value = 2670.46398
$ws.Range("N131").Value = -12750.46398

# CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 5240.2085
$ws.Range("I133").Value = 3875.8333
$ws.Range("J133").Value = 9333.333000000001
$ws.Range("K133").Value = 11627.4999
$ws.Range("L133").Value = 27999.999
$ws.Range("M133").Value = -6567.499899999999
$ws.Range("N133").Value = -38119.999

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 753.64703
$ws.Range("J135").Value = 1400
$ws.Range("L135").Value = 12600
$ws.Range("N135").Value = -17670

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 16919784
$ws.Range("I137").Value = 1997.9333
$ws.Range("J137").Value = 27070458
$ws.Range("K137").Value = 5993.7999
$ws.Range("L137").Value = 81211374
$ws.Range("M137").Value = -893.7999
$ws.Range("N137").Value = -81221574

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2948.32
$ws.Range("I140").Value = 907.7143
$ws.Range("K140").Value = 2723.1429
$ws.Range("M140").Value = 2456.8571

# LTW row 50
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = 0

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 44409.562
$ws.Range("I93").Value = 640.9091
$ws.Range("J93").Value = 140700.6
$ws.Range("K93").Value = 640.9091
$ws.Range("L93").Value = 140700.6
$ws.Range("M93").Value = 607.0909
$ws.Range("N93").Value = -143196.6

# WVR row 61
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 13346
$ws.Range("J61").Value = 21903.5
$ws.Range("L61").Value = 21903.5
$ws.Range("N61").Value = -22487.5

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4000
$ws.Range("L62").ClearContents()
$ws.Range("N62").Value = 0
$ws.Range("M62").Value = -3376

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 20000
$ws.Range("L65").ClearContents()
$ws.Range("N65").Value = 0
$ws.Range("M65").Value = -16880

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5297.476
$ws.Range("I107").Value = 6840.4375
$ws.Range("J107").Value = 360
$ws.Range("K107").Value = 20521.3125
$ws.Range("L107").Value = 1080
$ws.Range("M107").Value = -18601.3125
$ws.Range("N107").Value = -4920
